$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update puzzle numbers ---
$ws.Range("G3").Value = 5
$ws.Range("K3").Value = 7
$ws.Range("I5").Value = 6
$ws.Range("K5").Value = 9
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 8
$ws.Range("K7").Value = 1

# --- Update selection ---
$ws.Range("I16").Select()
